$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Neg_Change")
$ws2 = $wb.Worksheets.Item("Pos_Change")

# --- Sheet1 (Neg_Change): replace data rows 2-9 ---
$ws1.Range("A2").Value = "MARUTI"
$ws1.Range("B2").Value = 12840
$ws1.Range("C2").Value = 12925
$ws1.Range("D2").Value = 12752
$ws1.Range("E2").Value = 12839
$ws1.Range("F2").Value = 252065
$ws1.Range("G2").Value = 504476
$ws1.Range("H2").Value = -0.5003429300898358
$ws1.Range("I2").Value = "MARUTI"
$ws1.Range("A3").Value = "SIEMENS"
$ws1.Range("B3").Value = 3136.9
$ws1.Range("C3").Value = 3192
$ws1.Range("D3").Value = 3102.2
$ws1.Range("E3").Value = 3150
$ws1.Range("F3").Value = 305836
$ws1.Range("G3").Value = 697711
$ws1.Range("H3").Value = -0.5616580503962242
$ws1.Range("I3").Value = "SIEMENS"
$ws1.Range("A4").Value = "DLF"
$ws1.Range("B4").Value = 758.5
$ws1.Range("C4").Value = 767.2
$ws1.Range("D4").Value = 753.55
$ws1.Range("E4").Value = 757.7
$ws1.Range("F4").Value = 1366582
$ws1.Range("G4").Value = 2771092
$ws1.Range("H4").Value = -0.5068435115109856
$ws1.Range("I4").Value = "DLF"
$ws1.Range("A5").Value = "AUBANK"
$ws1.Range("B5").Value = 730
$ws1.Range("C5").Value = 742.5
$ws1.Range("D5").Value = 728.95
$ws1.Range("E5").Value = 740
$ws1.Range("F5").Value = 1328439
$ws1.Range("G5").Value = 2828419
$ws1.Range("H5").Value = -0.5303245381960735
$ws1.Range("I5").Value = "AUBANK"
$ws1.Range("A6").Value = "POLYCAB"
$ws1.Range("B6").Value = 6827
$ws1.Range("C6").Value = 6980.5
$ws1.Range("D6").Value = 6820
$ws1.Range("E6").Value = 6942
$ws1.Range("F6").Value = 255672
$ws1.Range("G6").Value = 522463
$ws1.Range("H6").Value = -0.5106409449090175
$ws1.Range("I6").Value = "POLYCAB"
$ws1.Range("A7").Value = "TIINDIA"
$ws1.Range("B7").Value = 3021.8
$ws1.Range("C7").Value = 3071.9
$ws1.Range("D7").Value = 3021.8
$ws1.Range("E7").Value = 3052
$ws1.Range("F7").Value = 187828
$ws1.Range("G7").Value = 469315
$ws1.Range("H7").Value = -0.5997826619647785
$ws1.Range("I7").Value = "TIINDIA"
$ws1.Range("A8").Value = "PERSISTENT"
$ws1.Range("B8").Value = 5275
$ws1.Range("C8").Value = 5277.5
$ws1.Range("D8").Value = 5202
$ws1.Range("E8").Value = 5240
$ws1.Range("F8").Value = 242009
$ws1.Range("G8").Value = 553305
$ws1.Range("H8").Value = -0.5626119409728811
$ws1.Range("I8").Value = "PERSISTENT"
$ws1.Range("A9").Value = "INOXWIND"
$ws1.Range("B9").Value = 141.45
$ws1.Range("C9").Value = 142.08
$ws1.Range("D9").Value = 137.7
$ws1.Range("E9").Value = 138.3
$ws1.Range("F9").Value = 4652786
$ws1.Range("G9").Value = 9769604
$ws1.Range("H9").Value = -0.5237487619764322
$ws1.Range("I9").Value = "INOXWIND"

# --- Sheet2 (Pos_Change): replace data rows 2-20 ---
$ws2.Range("A2").Value = "BHARTIARTL"
$ws2.Range("B2").Value = 1853
$ws2.Range("C2").Value = 1872.5
$ws2.Range("D2").Value = 1851.8
$ws2.Range("E2").Value = 1870.9
$ws2.Range("F2").Value = 6697206
$ws2.Range("G2").Value = 4645011
$ws2.Range("H2").Value = 0.441806273440472
$ws2.Range("I2").Value = "BHARTIARTL"
$ws2.Range("A3").Value = "ICICIBANK"
$ws2.Range("B3").Value = 1424.6
$ws2.Range("C3").Value = 1427.7
$ws2.Range("D3").Value = 1415.2
$ws2.Range("E3").Value = 1421.8
$ws2.Range("F3").Value = 9009465
$ws2.Range("G3").Value = 5913947
$ws2.Range("H3").Value = 0.5234267402125856
$ws2.Range("I3").Value = "ICICIBANK"
$ws2.Range("A4").Value = "INDHOTEL"
$ws2.Range("B4").Value = 751.1
$ws2.Range("C4").Value = 771.3
$ws2.Range("D4").Value = 750
$ws2.Range("E4").Value = 770.25
$ws2.Range("F4").Value = 2637612
$ws2.Range("G4").Value = 1747197
$ws2.Range("H4").Value = 0.5096248448228792
$ws2.Range("I4").Value = "INDHOTEL"
$ws2.Range("A5").Value = "LODHA"
$ws2.Range("B5").Value = 1214.2
$ws2.Range("C5").Value = 1230.7
$ws2.Range("D5").Value = 1214.2
$ws2.Range("E5").Value = 1229.5
$ws2.Range("F5").Value = 964805
$ws2.Range("G5").Value = 673722
$ws2.Range("H5").Value = 0.4320520927029249
$ws2.Range("I5").Value = "LODHA"
$ws2.Range("A6").Value = "TATAPOWER"
$ws2.Range("B6").Value = 385.5
$ws2.Range("C6").Value = 388.95
$ws2.Range("D6").Value = 384.4
$ws2.Range("E6").Value = 387.1
$ws2.Range("F6").Value = 3757839
$ws2.Range("G6").Value = 2575426
$ws2.Range("H6").Value = 0.4591135602420726
$ws2.Range("I6").Value = "TATAPOWER"
$ws2.Range("A7").Value = "ICICIGI"
$ws2.Range("B7").Value = 1902.9
$ws2.Range("C7").Value = 1924.3
$ws2.Range("D7").Value = 1901
$ws2.Range("E7").Value = 1913.5
$ws2.Range("F7").Value = 335318
$ws2.Range("G7").Value = 233350
$ws2.Range("H7").Value = 0.4369745018212985
$ws2.Range("I7").Value = "ICICIGI"
$ws2.Range("A8").Value = "LTIM"
$ws2.Range("B8").Value = 5130
$ws2.Range("C8").Value = 5132
$ws2.Range("D8").Value = 5055.5
$ws2.Range("E8").Value = 5099.5
$ws2.Range("F8").Value = 237765
$ws2.Range("G8").Value = 156885
$ws2.Range("H8").Value = 0.5155368582082417
$ws2.Range("I8").Value = "LTIM"
$ws2.Range("A9").Value = "VOLTAS"
$ws2.Range("B9").Value = 1249.9
$ws2.Range("C9").Value = 1284
$ws2.Range("D9").Value = 1243.7
$ws2.Range("E9").Value = 1275.7
$ws2.Range("F9").Value = 1199406
$ws2.Range("G9").Value = 839568
$ws2.Range("H9").Value = 0.4285989937682236
$ws2.Range("I9").Value = "VOLTAS"
$ws2.Range("A10").Value = "SONACOMS"
$ws2.Range("B10").Value = 442.55
$ws2.Range("C10").Value = 450.55
$ws2.Range("D10").Value = 442.05
$ws2.Range("E10").Value = 449.1
$ws2.Range("F10").Value = 1320587
$ws2.Range("G10").Value = 921533
$ws2.Range("H10").Value = 0.4330327834163291
$ws2.Range("I10").Value = "SONACOMS"
$ws2.Range("A11").Value = "POLICYBZR"
$ws2.Range("B11").Value = 1835.1
$ws2.Range("C11").Value = 1889.5
$ws2.Range("D11").Value = 1825.5
$ws2.Range("E11").Value = 1852.5
$ws2.Range("F11").Value = 1329799
$ws2.Range("G11").Value = 875077
$ws2.Range("H11").Value = 0.5196365576972084
$ws2.Range("I11").Value = "POLICYBZR"
$ws2.Range("A12").Value = "DIXON"
$ws2.Range("B12").Value = 15900
$ws2.Range("C12").Value = 16070
$ws2.Range("D12").Value = 15593
$ws2.Range("E12").Value = 15980
$ws2.Range("F12").Value = 235815
$ws2.Range("G12").Value = 164396
$ws2.Range("H12").Value = 0.4344327112581814
$ws2.Range("I12").Value = "DIXON"
$ws2.Range("A13").Value = "MARICO"
$ws2.Range("B13").Value = 705.05
$ws2.Range("C13").Value = 711
$ws2.Range("D13").Value = 699.25
$ws2.Range("E13").Value = 710.25
$ws2.Range("F13").Value = 1291837
$ws2.Range("G13").Value = 867153
$ws2.Range("H13").Value = 0.4897451776099489
$ws2.Range("I13").Value = "MARICO"
$ws2.Range("A14").Value = "LICHSGFIN"
$ws2.Range("B14").Value = 576.1
$ws2.Range("C14").Value = 577
$ws2.Range("D14").Value = 571.35
$ws2.Range("E14").Value = 574.8
$ws2.Range("F14").Value = 939214
$ws2.Range("G14").Value = 616119
$ws2.Range("H14").Value = 0.5244035648957426
$ws2.Range("I14").Value = "LICHSGFIN"
$ws2.Range("A15").Value = "PAGEIND"
$ws2.Range("B15").Value = 44010
$ws2.Range("C15").Value = 44280
$ws2.Range("D15").Value = 43405
$ws2.Range("E15").Value = 43600
$ws2.Range("F15").Value = 36102
$ws2.Range("G15").Value = 24274
$ws2.Range("H15").Value = 0.4872703303946609
$ws2.Range("I15").Value = "PAGEIND"
$ws2.Range("A16").Value = "IDEA"
$ws2.Range("B16").Value = 6.53
$ws2.Range("C16").Value = 6.54
$ws2.Range("D16").Value = 6.35
$ws2.Range("E16").Value = 6.38
$ws2.Range("F16").Value = 470360407
$ws2.Range("G16").Value = 317847247
$ws2.Range("H16").Value = 0.4798316217601218
$ws2.Range("I16").Value = "IDEA"
$ws2.Range("A17").Value = "LAURUSLABS"
$ws2.Range("B17").Value = 836.6
$ws2.Range("C17").Value = 866
$ws2.Range("D17").Value = 833.05
$ws2.Range("E17").Value = 866
$ws2.Range("F17").Value = 2817582
$ws2.Range("G17").Value = 1863616
$ws2.Range("H17").Value = 0.5118897884542738
$ws2.Range("I17").Value = "LAURUSLABS"
$ws2.Range("A18").Value = "IEX"
$ws2.Range("B18").Value = 138
$ws2.Range("C18").Value = 142.29
$ws2.Range("D18").Value = 137.81
$ws2.Range("E18").Value = 141.1
$ws2.Range("F18").Value = 16615911
$ws2.Range("G18").Value = 10917993
$ws2.Range("H18").Value = 0.5218832801962778
$ws2.Range("I18").Value = "IEX"
$ws2.Range("A19").Value = "NUVAMA"
$ws2.Range("B19").Value = 6929
$ws2.Range("C19").Value = 7000
$ws2.Range("D19").Value = 6720
$ws2.Range("E19").Value = 6955
$ws2.Range("F19").Value = 133141
$ws2.Range("G19").Value = 85228
$ws2.Range("H19").Value = 0.5621744027784297
$ws2.Range("I19").Value = "NUVAMA"
$ws2.Range("A20").Value = "NCC"
$ws2.Range("B20").Value = 222.6
$ws2.Range("C20").Value = 225.09
$ws2.Range("D20").Value = 221.1
$ws2.Range("E20").Value = 222.27
$ws2.Range("F20").Value = 2001471
$ws2.Range("G20").Value = 1278810
$ws2.Range("H20").Value = 0.5651042766321814
$ws2.Range("I20").Value = "NCC"

Write-Host "Update complete."
